$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.275184154510498
$ws.Range("B1").Value = 2.13282299041748
$ws.Range("C1").Value = 4.73020076751709
$ws.Range("D1").Value = 3.31791353225708
$ws.Range("E1").Value = 1.375092148780823
